$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 230, shifting the existing rows 230..335 down to 231..336.
$ws.Rows(230).Insert()

# Populate the newly inserted row 230 with the new weekly record.
$ws.Cells.Item(230, 1).Value = 3
$ws.Cells.Item(230, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(230, 3).Value = "Coquimbo"
$ws.Cells.Item(230, 4).Value = 44609
$ws.Cells.Item(230, 5).Value = 5
$ws.Cells.Item(230, 6).Value = 100112017
$ws.Cells.Item(230, 7).Value = "Apio"
$ws.Cells.Item(230, 8).Value = "Americana (o)"
$ws.Cells.Item(230, 9).Value = "Primera"
$ws.Cells.Item(230, 10).Value = 230
$ws.Cells.Item(230, 11).Value = 9000
$ws.Cells.Item(230, 12).Value = 9500
$ws.Cells.Item(230, 13).Value = 9239
$ws.Cells.Item(230, 14).Value = "$/docena de matas"
$ws.Cells.Item(230, 15).Value = "Pan de Az$([char]0xFA)car"
$ws.Cells.Item(230, 16).Value = 1540
$ws.Cells.Item(230, 17).Value = 6
$ws.Cells.Item(230, 18).Value = "Hortaliza"

# Keep the date formatted the same way as the rest of the Fecha column.
$ws.Cells.Item(230, 4).NumberFormat = $ws.Cells.Item(231, 4).NumberFormat
